# Slide 1 ("CSSE 220 Day 10"): the "Today's Attendance password" textbox
# currently shows a blank-line placeholder ("__________") under the
# prompt. Fill in today's attendance password ("uml") in its place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$passwordBox = $s.Shapes.Item("TextBox 1")
$tr = $passwordBox.TextFrame.TextRange

# The textbox has two paragraphs:
#   1) "Today's Attendance password"
#   2) "__________"
# Locate the second paragraph's run (the blank-line placeholder) and
# replace just its text, leaving the rest of the run formatting
# (highlight, size, etc.) untouched.
$fullText = $tr.Text
$placeholder = "__________"
$startIdx = $fullText.IndexOf($placeholder)

if ($startIdx -ge 0) {
    $target = $tr.Characters($startIdx + 1, $placeholder.Length)
    $target.Text = "uml"
}
